# Scheduled runner refresh: update Leve profit-calc columns (H:N)
# with latest market-board price snapshots across several sheets.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# row 88
$ws.Range("H88").Value = 3710.2083
$ws.Range("I88").Value = 2750
$ws.Range("J88").Value = 3751.9565
$ws.Range("K88").Value = 2750
$ws.Range("L88").Value = 3751.9565
$ws.Range("M88").Value = -2344
$ws.Range("N88").Value = -4563.9565
# row 91
$ws.Range("H91").Value = 3710.2083
$ws.Range("I91").Value = 2750
$ws.Range("J91").Value = 3751.9565
$ws.Range("K91").Value = 2750
$ws.Range("L91").Value = 3751.9565
$ws.Range("M91").Value = -1346
$ws.Range("N91").Value = -6559.9565
# row 98
$ws.Range("H98").Value = 1824.0714
$ws.Range("I98").Value = 1887.8846
$ws.Range("J98").Value = 994.5
$ws.Range("K98").Value = 1887.8846
$ws.Range("L98").Value = 994.5
$ws.Range("M98").Value = -389.8846000000001
$ws.Range("N98").Value = -3990.5
# row 112
$ws.Range("H112").Value = 4002.3416
$ws.Range("J112").Value = 4077.4
$ws.Range("L112").Value = 12232.2
$ws.Range("N112").Value = -14448.2
# row 122
$ws.Range("H122").Value = 1824.0714
$ws.Range("I122").Value = 1887.8846
$ws.Range("J122").Value = 994.5
$ws.Range("K122").Value = 5663.6538
$ws.Range("L122").Value = 2983.5
$ws.Range("M122").Value = -3213.6538
$ws.Range("N122").Value = -7883.5
# row 132
$ws.Range("H132").Value = 11908149
$ws.Range("I132").Value = 13161338
$ws.Range("K132").Value = 39484014
$ws.Range("M132").Value = -39481484
# row 133
$ws.Range("H133").Value = 122847.25
$ws.Range("J133").Value = 122847.25
$ws.Range("L133").Value = 122847.25
$ws.Range("N133").Value = -132967.25
# row 134
$ws.Range("H134").Value = 112280.414
$ws.Range("J134").Value = 112280.414
$ws.Range("L134").Value = 112280.414
$ws.Range("N134").Value = -122420.414
# row 137
$ws.Range("H137").Value = 34569.37
$ws.Range("I137").Value = 41991.316
$ws.Range("J137").Value = 1912.8
$ws.Range("K137").Value = 125973.948
$ws.Range("L137").Value = 5738.4
$ws.Range("M137").Value = -123423.948
$ws.Range("N137").Value = -10838.4
# row 139
$ws.Range("H139").Value = 97047.836
$ws.Range("J139").Value = 100457.4
$ws.Range("L139").Value = 100457.4
$ws.Range("N139").Value = -110737.4

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 7689.1465
$ws.Range("I32").Value = 4509.4707
$ws.Range("J32").Value = 23133.285
$ws.Range("K32").Value = 4509.4707
$ws.Range("L32").Value = 23133.285
$ws.Range("M32").Value = -4222.4707
$ws.Range("N32").Value = -23707.285
# row 45
$ws.Range("H45").Value = 7996148
$ws.Range("I45").Value = 13079175
$ws.Range("J45").Value = 8533.857
$ws.Range("K45").Value = 13079175
$ws.Range("L45").Value = 8533.857
$ws.Range("M45").Value = -13078798
$ws.Range("N45").Value = -9287.857
# row 74
$ws.Range("H74").Value = 32057.861
$ws.Range("I74").Value = 17276.215
$ws.Range("K74").Value = 17276.215
$ws.Range("M74").Value = -16402.215
# row 77
$ws.Range("H77").Value = 32057.861
$ws.Range("I77").Value = 17276.215
$ws.Range("K77").Value = 86381.075
$ws.Range("M77").Value = -82013.075
# row 140
$ws.Range("H140").Value = 117481.336
$ws.Range("J140").Value = 117481.336
$ws.Range("L140").Value = 117481.336
$ws.Range("N140").Value = -127841.336

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# row 64
$ws.Range("H64").Value = 2299.6667
$ws.Range("I64").Value = 2199
$ws.Range("J64").Value = 2350
$ws.Range("K64").Value = 2199
$ws.Range("L64").Value = 2350
$ws.Range("M64").Value = -1974
$ws.Range("N64").Value = -2800
# row 67
$ws.Range("H67").Value = 2299.6667
$ws.Range("I67").Value = 2199
$ws.Range("J67").Value = 2350
$ws.Range("K67").Value = 2199
$ws.Range("L67").Value = 2350
$ws.Range("M67").Value = -1419
$ws.Range("N67").Value = -3910
# row 140
$ws.Range("H140").Value = 76173
$ws.Range("J140").Value = 76173
$ws.Range("L140").Value = 76173
$ws.Range("N140").Value = -86533

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 17356.953
$ws.Range("I31").Value = 2363.2856
$ws.Range("J31").Value = 47344.285
$ws.Range("K31").Value = 2363.2856
$ws.Range("L31").Value = 47344.285
$ws.Range("M31").Value = -2068.2856
$ws.Range("N31").Value = -47934.285
# row 34
$ws.Range("H34").Value = 17356.953
$ws.Range("I34").Value = 2363.2856
$ws.Range("J34").Value = 47344.285
$ws.Range("K34").Value = 2363.2856
$ws.Range("L34").Value = 47344.285
$ws.Range("M34").Value = -2161.2856
$ws.Range("N34").Value = -47748.285
# row 58
$ws.Range("H58").Value = 2239.0962
$ws.Range("I58").Value = 2175.8
$ws.Range("J58").Value = 2369.4119
$ws.Range("K58").Value = 2175.8
$ws.Range("L58").Value = 2369.4119
$ws.Range("M58").Value = -1972.8
$ws.Range("N58").Value = -2775.4119
# row 132
$ws.Range("H132").Value = 40800.04
$ws.Range("I132").Value = 24202.303
$ws.Range("J132").Value = 112170.3
$ws.Range("K132").Value = 72606.909
$ws.Range("L132").Value = 336510.9
$ws.Range("M132").Value = -70076.909
$ws.Range("N132").Value = -341570.9
# row 135
$ws.Range("H135").Value = 148374
$ws.Range("J135").Value = 148374
$ws.Range("L135").Value = 148374
$ws.Range("N135").Value = -158514
# row 136
$ws.Range("H136").Value = 2239.0962
$ws.Range("I136").Value = 2175.8
$ws.Range("J136").Value = 2369.4119
$ws.Range("K136").Value = 6527.400000000001
$ws.Range("L136").Value = 7108.2357
$ws.Range("M136").Value = -3977.400000000001
$ws.Range("N136").Value = -12208.2357
# row 140
$ws.Range("H140").Value = 57722.5
$ws.Range("J140").Value = 57722.5
$ws.Range("L140").Value = 57722.5
$ws.Range("N140").Value = -68082.5
# row 141
$ws.Range("H141").Value = 43560.8
$ws.Range("J141").Value = 43560.8
$ws.Range("L141").Value = 43560.8
$ws.Range("N141").Value = -53920.8

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# row 132
$ws.Range("H132").Value = 2397.7646
$ws.Range("I132").Value = 2544.6206
$ws.Range("J132").Value = 1546
$ws.Range("K132").Value = 7633.861800000001
$ws.Range("L132").Value = 4638
$ws.Range("M132").Value = -5103.861800000001
$ws.Range("N132").Value = -9698
# row 139
$ws.Range("H139").Value = 56980.125
$ws.Range("J139").Value = 56980.125
$ws.Range("L139").Value = 56980.125
$ws.Range("N139").Value = -67260.125

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# row 136
$ws.Range("H136").Value = 40023.094
$ws.Range("I136").Value = 56598.473
$ws.Range("J136").Value = 4922.294
$ws.Range("K136").Value = 169795.419
$ws.Range("L136").Value = 14766.882
$ws.Range("M136").Value = -167245.419
$ws.Range("N136").Value = -19866.882
# row 139
$ws.Range("H139").Value = 65988
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()
# row 140
$ws.Range("H140").Value = 93844.86
$ws.Range("J140").Value = 93137.336
$ws.Range("L140").Value = 93137.336
$ws.Range("N140").Value = -103497.336
# row 141
$ws.Range("H141").Value = 112121.07
$ws.Range("J141").Value = 112121.07
$ws.Range("L141").Value = 112121.07
$ws.Range("N141").Value = -122481.07

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# row 119
$ws.Range("H119").Value = 29000
$ws.Range("J119").Value = 29000
$ws.Range("L119").Value = 29000
$ws.Range("N119").Value = -38676
# row 132
$ws.Range("H132").Value = 24070400
$ws.Range("I132").Value = 33334400
$ws.Range("J132").Value = 910401.8
$ws.Range("K132").Value = 100003200
$ws.Range("L132").Value = 2731205.4
$ws.Range("M132").Value = -100000670
$ws.Range("N132").Value = -2736265.4
# row 136
$ws.Range("H136").Value = 2179.4834
$ws.Range("J136").Value = 3398.75
$ws.Range("L136").Value = 10196.25
$ws.Range("N136").Value = -15296.25
# row 141
$ws.Range("H141").Value = 131666.33
$ws.Range("J141").Value = 131666.33
$ws.Range("L141").Value = 131666.33
$ws.Range("N141").Value = -142026.33
